$wb = $excel.ActiveWorkbook

# 1. Move the "Fuel_type" sheet so it sits right before "Fuel Pump"
#    (new order: User, Shift, Fuel_type, Fuel Pump, Price, Transaction,
#     Fuel Tank Inventory, Deliveries)
#    NOTE: worksheet object handles here are index-bound, so after a
#    Move() re-fetch sheets by name rather than reusing old variables.
$wsFuelType = $wb.Worksheets.Item("Fuel_type")
$wsFuelPump = $wb.Worksheets.Item("Fuel Pump")
$wsFuelType.Move($wsFuelPump)

# 2. Remove the "Volume_delivered" column (E) from "Fuel Tank Inventory"
$wsTank = $wb.Worksheets.Item("Fuel Tank Inventory")
$wsTank.Range("E1:E4").ClearContents()

# 3. Add the new header row to "Deliveries"
$wsDeliveries = $wb.Worksheets.Item("Deliveries")
$wsDeliveries.Range("A1").Value = "Delivery_id"
$wsDeliveries.Range("B1").Value = "Fuel_type_id"
$wsDeliveries.Range("C1").Value = "Volume"
$wsDeliveries.Range("D1").Value = "Supplier"
$wsDeliveries.Columns.Item(1).ColumnWidth = 11.333333333333332
$wsDeliveries.Columns.Item(2).ColumnWidth = 11.5
$wsDeliveries.Columns.Item(3).ColumnWidth = 8.666666666666666

# 4. Restore / adjust each sheet's selection (active cell) to match the
#    target state.
$wsShift = $wb.Worksheets.Item("Shift")
$wsShift.Range("G16").Select()

$wsFuelTypeAfter = $wb.Worksheets.Item("Fuel_type")
$wsFuelTypeAfter.Range("C1").Select()

$wsTank.Range("E1").Select()

$wsDeliveries.Range("F7").Select()

# 5. Make "Fuel Pump" the active sheet/tab with its target selection -
#    doing this last makes it the active tab in the saved workbook.
$wsFuelPumpAfter = $wb.Worksheets.Item("Fuel Pump")
$wsFuelPumpAfter.Activate()
$wsFuelPumpAfter.Range("H5").Select()
